$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Formula = "=IF(C2<=50000,0,IF(C2<=100000,(C2-50000)*10%,50000*0.1+(C2-100000)*20%))"
$ws.Range("D3:D10").Formula = "=IF(C3<=50000,0,IF(C3<=100000,(C3-50000)*10%,50000*0.1+(C3-100000)*20%))"

$ws.Range("C13").Formula = '=COUNTIFS(B2:B10,1,C2:C10,">50000",C2:C10,"<=100000")'
